$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "003"
$ws.Range("N2").Value = "2017-03-31 00:00:00"
$ws.Range("O2").Value = 13188261.22
$ws.Range("P2").Value = 155799164.53
$ws.Range("Q2").Value = 139740128.37
$ws.Range("R2").Value = 79.1724692623
$ws.Range("S2").Value = 103067219.91
$ws.Range("T2").Value = 103067219.91
$ws.Range("U2").Value = 60.4653990737
$ws.Range("V2").Value = 7797695.17
$ws.Range("W2").Value = 22601506.34
$ws.Range("X2").Value = 5961936.11
$ws.Range("Y2").Value = 16059036.16
$ws.Range("Z2").Value = 17419226.89
$ws.Range("AA2").Value = 4303285.44
$ws.Range("AG2").Value = 311770.84
$ws.Range("AP2").Value = 68.9635825191
$ws.Range("AQ2").Value = 12.958478971384
$ws.Range("AR2").Value = 15.790294826551
$ws.Range("AS2").Value = 12167455.89
$ws.Range("AT2").Value = 6.790128099232
